# Shipper Contact Code Commit
# Adds Expected/Actual Result + Status columns to the "Carrier Details" sheet,
# writes a success message, and updates sheet selections.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Carrier Details")
$ws2 = $wb.Worksheets.Item("CustomizeGrid")
$ws3 = $wb.Worksheets.Item("View Carrier Details")

# --- Sheet1: "Carrier Details" -------------------------------------------
# New header cells (H1:J1) mirroring the "View Carrier Details" headers.
$ws1.Range("H1").Value = "Expected Result"
$ws1.Range("I1").Value = "Actual Result"
$ws1.Range("J1").Value = "Status"

# Bold-bordered header style (matches existing header formatting elsewhere).
$headerRange = $ws1.Range("H1:J1")
$headerRange.NumberFormat = "@"
$headerRange.Borders.LineStyle = 1

# Result message cell.
$ws1.Range("H2").Value = "Company Added Successfully"
$ws1.Range("H2").NumberFormat = "@"
$ws1.Range("H2").Borders(7).LineStyle = 1
$ws1.Range("H2").Borders(10).LineStyle = 1

# Column widths for the newly populated columns.
$ws1.Columns.Item(8).ColumnWidth = 27.28515625
$ws1.Columns.Item(9).ColumnWidth = 12.5703125

$ws1.Range("I5").Select()

# --- Sheet2: "CustomizeGrid" ----------------------------------------------
$ws2.Range("A1:F2").Select()

# --- Sheet3: "View Carrier Details" ---------------------------------------
$ws3.Range("A1:J2").Select()
